$data = New-Object 'object[,]' 39,2
$data[0,0] = "label"
$data[0,1] = "path"
$data[1,0] = "CANDIOLO_HG19"
$data[1,1] = "data/source_data/TCGA_RNAseq_hg19/CANDIOLO_HG19.rds"
$data[2,0] = "GMQL_GRCH38"
$data[2,1] = "data/source_data/TCGA_RNAseq_GRCh38/TCGA_RNAseq_GRCh38.rds"
$data[3,0] = "GMQL_GRCH38_FILT"
$data[3,1] = "data/source_data/TCGA_RNAseq_GRCh38/GMQL_GRCH38_FILTERED.rds"
$data[4,0] = "GMQL_GRCH38_FILT_UNIF"
$data[4,1] = "data/source_data/TCGA_RNAseq_GRCh38/GMQL_GRCH38_FILTERED_unif.rds"
$data[5,0] = "GMQL_GRCH38_META"
$data[5,1] = "data/source_data/TCGA_RNAseq_GRCh38/annotations/gmql_grch38_metadata.xlsx"
$data[6,0] = "GMQL_GRCH38_ANNOT"
$data[6,1] = "data/source_data/TCGA_RNAseq_GRCh38/annotations/gmql_grch38_annotations.xlsx"
$data[7,0] = "PDX_1"
$data[7,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod1_LMX.rds"
$data[8,0] = "PDX_2"
$data[8,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod2_LMX.rds"
$data[9,0] = "PDX_3"
$data[9,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod3_LMX.rds"
$data[10,0] = "PDX_4"
$data[10,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod4_LMX.rds"
$data[11,0] = "PDX_5"
$data[11,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod5_LMX.rds"
$data[12,0] = "PDX_6"
$data[12,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod6_LMX.rds"
$data[13,0] = "PDX_1_FILTERED"
$data[13,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod1_LMX_FILTERED.rds"
$data[14,0] = "PDX_2_FILTERED"
$data[14,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod2_LMX_FILTERED.rds"
$data[15,0] = "PDX_3_FILTERED"
$data[15,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod3_LMX_FILTERED.rds"
$data[16,0] = "PDX_4_FILTERED"
$data[16,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod4_LMX_FILTERED.rds"
$data[17,0] = "PDX_5_FILTERED"
$data[17,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod5_LMX.rds"
$data[18,0] = "PDX_6_FILTERED"
$data[18,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod6_LMX.rds"
$data[19,0] = "PDX_MERGED_FILT"
$data[19,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod_MERGED_FILTERED.rds"
$data[20,0] = "PDX_MERGED_FILT_UNIF"
$data[20,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/Hbiod_MERGED_FILTERED_unif.rds"
$data[21,0] = "PDX_MERGED_ANNOT"
$data[21,1] = "data/source_data/PDX_RNAseq_GRCh38/PDX_GRCH38/annotations/pdx_grch38_annotations.xlsx"
$data[22,0] = "NTP_REF_TCGA"
$data[22,1] = "data/references/NTP_cpm_tcga_reference.rds"
$data[23,0] = "NTP_REF_PDX"
$data[23,1] = "data/references/NTP_cpm_pdx_reference.rds"
$data[24,0] = "TCGA_SPLITTING"
$data[24,1] = "data/references/tcga_splitting.xlsx"
$data[25,0] = "FEATURES_ORIGINAL"
$data[25,1] = "data/genes/ntp_signature_tcga_hg19_original.xlsx"
$data[26,0] = "FEATURES_HG19"
$data[26,1] = "data/genes/ntp_signature_tcga_hg19_filtered.xlsx"
$data[27,0] = "FEATURES_GRCH38"
$data[27,1] = "data/genes/ntp_signature_tcga_grch38.xlsx"
$data[28,0] = "FEATURES_PDX"
$data[28,1] = "data/genes/ntp_signature_pdx_grch38.xlsx"
$data[29,0] = "PUB_NTP"
$data[29,1] = "data/references/published_ntp.xlsx"
$data[30,0] = "PUB_TSP"
$data[30,1] = "data/references/published_tsp.xlsx"
$data[31,0] = "NTP_THR"
$data[31,1] = "data/references/ml_ntp_thresholds.rds"
$data[32,0] = "HD_PDX"
$data[32,1] = "data/source_data/PDX_RNAseq_GRCh38/HIGH_DEPTH_PDX_GRCH38/high_depth_LMX_samples.xlsx"
$data[33,0] = "BIO_LASSO_TCGA"
$data[33,1] = "data/genes/feature_selection/bio_driven_lasso_tcga_robust_union.xlsx"
$data[34,0] = "BIO_LASSO_PDX"
$data[34,1] = "data/genes/feature_selection/bio_driven_lasso_pdx_robust_union.xlsx"
$data[35,0] = "BIO_DRIVEN_TCGA"
$data[35,1] = "data/genes/feature_selection/bio_fs_genes_tcga.xlsx"
$data[36,0] = "BIO_DRIVEN_PDX"
$data[36,1] = "data/genes/feature_selection/bio_fs_genes_pdx.xlsx"
$data[37,0] = "ENTREZ_CORR_TCGA"
$data[37,1] = "data/genes/entrez_corr_tcga.xlsx"
$data[38,0] = "ENTREZ_CORR_PDX"
$data[38,1] = "data/genes/entrez_corr_pdx.xlsx"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the full label/path table (header + 38 data rows)
$ws.Range("A1:B39").Value = $data

# Column A width changed from 44.33203125 to 27.6640625 (characters, stored width).
# The COM ColumnWidth setter only offers coarse-grained pixel snapping in this
# runtime, so we pick the ColumnWidth value whose resulting stored width is the
# closest achievable approximation to the target.
$ws.Columns.Item(1).ColumnWidth = 26.83

# Restore the selection / active cell as last edited by the author
$ws.Range("B34").Select()
